# Update "想去人数" (want-to-go count) figures in column F across all sheets.
# This mirrors a refreshed data pull (gh-pages output regenerated at 802b57d).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 11394
$ws1.Range("F3").Value  = 1910
$ws1.Range("F7").Value  = 753
$ws1.Range("F8").Value  = 995
$ws1.Range("F9").Value  = 584
$ws1.Range("F10").Value = 447
$ws1.Range("F11").Value = 1320
$ws1.Range("F12").Value = 661
$ws1.Range("F13").Value = 104
$ws1.Range("F14").Value = 7
$ws1.Range("F15").Value = 972
$ws1.Range("F16").Value = 518
$ws1.Range("F17").Value = 664
$ws1.Range("F18").Value = 1078
$ws1.Range("F19").Value = 204
$ws1.Range("F20").Value = 927
$ws1.Range("F21").Value = 131
$ws1.Range("F22").Value = 268
$ws1.Range("F24").Value = 258
$ws1.Range("F25").Value = 462
$ws1.Range("F26").Value = 491
$ws1.Range("F27").Value = 668
$ws1.Range("F28").Value = 172
$ws1.Range("F29").Value = 104
$ws1.Range("F30").Value = 320

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 36
$ws2.Range("F5").Value  = 879
$ws2.Range("F7").Value  = 5
$ws2.Range("F10").Value = 357

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 22

# Sheet "全部类型" (All types) - combined listing
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 11394
$ws4.Range("F3").Value  = 1910
$ws4.Range("F8").Value  = 753
$ws4.Range("F9").Value  = 995
$ws4.Range("F10").Value = 36
$ws4.Range("F11").Value = 584
$ws4.Range("F12").Value = 447
$ws4.Range("F13").Value = 22
$ws4.Range("F14").Value = 1320
$ws4.Range("F16").Value = 661
$ws4.Range("F17").Value = 104
$ws4.Range("F18").Value = 879
$ws4.Range("F19").Value = 7
$ws4.Range("F20").Value = 972
$ws4.Range("F21").Value = 518
$ws4.Range("F22").Value = 664
$ws4.Range("F23").Value = 1078
$ws4.Range("F24").Value = 204
$ws4.Range("F25").Value = 927
$ws4.Range("F26").Value = 131
$ws4.Range("F27").Value = 268
$ws4.Range("F30").Value = 258
$ws4.Range("F31").Value = 5
$ws4.Range("F33").Value = 462
$ws4.Range("F34").Value = 491
$ws4.Range("F35").Value = 668
$ws4.Range("F36").Value = 172
$ws4.Range("F38").Value = 104
$ws4.Range("F39").Value = 357
$ws4.Range("F40").Value = 320
